$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New columns D:K get a uniform 9-character width (replaces the old
# best-fit width on column J alone).
$ws.Range("D1:K1").ColumnWidth = 8.1

# Extend the thin bottom border under row 3 into the new column K.
$ws.Range("J3").Copy()
$ws.Range("K3").PasteSpecial(-4122)

# Add the 2022 year header, copying J4's formatting (style 9).
$ws.Range("J4").Copy()
$ws.Range("K4").PasteSpecial(-4122)
$ws.Range("K4").Value = 2022

# Add the 2022 percentage value, copying J5's formatting (style 13).
$ws.Range("J5").Copy()
$ws.Range("K5").PasteSpecial(-4122)
$ws.Range("K5").Value = 0.11705180708279034

$ws.Range("J12").Select()
